# Apply edit: add a new worksheet "Sheet1" after the existing
# "Cohort_Retentio" sheet, make it the active/selected sheet, and
# populate it with a sparse-matrix style template of admission
# selection-factor data.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing worksheet so that it
# lands at the end of the sheet tab order (matches the target: sheets
# are Cohort_Retentio, Sheet1).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$new.Name = "Sheet1"

# Header row (row 1): columns A..T
$new.Cells.Item(1,1).Value = "Value"
$new.Cells.Item(1,2).Value = "rigor of secondary school record"
$new.Cells.Item(1,3).Value = "class rank"
$new.Cells.Item(1,4).Value = "academic GPA"
$new.Cells.Item(1,5).Value = "standardized test scores"
$new.Cells.Item(1,6).Value = "application Essay"
$new.Cells.Item(1,7).Value = "recommendation"
$new.Cells.Item(1,8).Value = "interview"
$new.Cells.Item(1,9).Value = "extracurricular activities"
$new.Cells.Item(1,10).Value = "talent/ability"
$new.Cells.Item(1,11).Value = "character/personal qualities"
$new.Cells.Item(1,12).Value = "first generation"
$new.Cells.Item(1,13).Value = "alumni/ae relation"
$new.Cells.Item(1,14).Value = "geographical residence"
$new.Cells.Item(1,15).Value = "state residency"
$new.Cells.Item(1,16).Value = "religious affiliation/commitment"
$new.Cells.Item(1,17).Value = "racial/ethnic status"
$new.Cells.Item(1,18).Value = "volunteer work"
$new.Cells.Item(1,19).Value = "work experience"
$new.Cells.Item(1,20).Value = "level of applicant’s interest"

# Row 2
$new.Cells.Item(2,1).Value = "Very Important"
$new.Cells.Item(2,2).Value = 1
$new.Cells.Item(2,3).Value = 0
$new.Cells.Item(2,4).Value = 0

# Row 3
$new.Cells.Item(3,1).Value = "Considered"
$new.Cells.Item(3,2).Value = 0
$new.Cells.Item(3,3).Value = 1
$new.Cells.Item(3,4).Value = 0

# Row 4
$new.Cells.Item(4,1).Value = "Very Important"
$new.Cells.Item(4,2).Value = 0
$new.Cells.Item(4,3).Value = 0
$new.Cells.Item(4,4).Value = 1
$new.Cells.Item(4,5).Value = 0
$new.Cells.Item(4,6).Value = 0
$new.Cells.Item(4,7).Value = 0
$new.Cells.Item(4,8).Value = 0
$new.Cells.Item(4,9).Value = 0
$new.Cells.Item(4,10).Value = 0
$new.Cells.Item(4,11).Value = 0
$new.Cells.Item(4,12).Value = 0
$new.Cells.Item(4,13).Value = 0
$new.Cells.Item(4,14).Value = 0
$new.Cells.Item(4,15).Value = 0
$new.Cells.Item(4,16).Value = 0
$new.Cells.Item(4,17).Value = 0
$new.Cells.Item(4,18).Value = 0
$new.Cells.Item(4,19).Value = 0
$new.Cells.Item(4,20).Value = 0

# Row 5
$new.Cells.Item(5,1).Value = "Considered"
$new.Cells.Item(5,5).Value = 1
$new.Cells.Item(5,6).Value = 0
$new.Cells.Item(5,7).Value = 0
$new.Cells.Item(5,8).Value = 0
$new.Cells.Item(5,9).Value = 0
$new.Cells.Item(5,10).Value = 0
$new.Cells.Item(5,11).Value = 0
$new.Cells.Item(5,12).Value = 0
$new.Cells.Item(5,13).Value = 0
$new.Cells.Item(5,14).Value = 0
$new.Cells.Item(5,15).Value = 0
$new.Cells.Item(5,16).Value = 0
$new.Cells.Item(5,17).Value = 0
$new.Cells.Item(5,18).Value = 0
$new.Cells.Item(5,19).Value = 0
$new.Cells.Item(5,20).Value = 0

# Select the full populated range and make this new sheet the active
# (selected) tab, matching the authored workbook state.
$null = $new.Range("A1:T5").Select()
$null = $new.Activate()
